$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1254.8334
$ws.Range("I39").Value = 1381.25
$ws.Range("J39").Value = 1002
$ws.Range("K39").Value = 4143.75
$ws.Range("L39").Value = 3006
$ws.Range("M39").Value = -3847.75
$ws.Range("N39").Value = -3598
$ws.Range("H53").Value = 393.27274
$ws.Range("I53").Value = 270.83334
$ws.Range("J53").Value = 540.2
$ws.Range("K53").Value = 270.83334
$ws.Range("L53").Value = 540.2
$ws.Range("M53").Value = 366.16666
$ws.Range("N53").Value = -1814.2
$ws.Range("H58").Value = 335
$ws.Range("I58").Value = 380.42856
$ws.Range("K58").Value = 1141.28568
$ws.Range("M58").Value = -991.28568
$ws.Range("H62").Value = 5491.2104
$ws.Range("I62").Value = 4214.5
$ws.Range("J62").Value = 6419.727
$ws.Range("K62").Value = 4214.5
$ws.Range("L62").Value = 6419.727
$ws.Range("M62").Value = -3590.5
$ws.Range("N62").Value = -7667.727
$ws.Range("H65").Value = 5491.2104
$ws.Range("I65").Value = 4214.5
$ws.Range("J65").Value = 6419.727
$ws.Range("K65").Value = 21072.5
$ws.Range("L65").Value = 32098.635
$ws.Range("M65").Value = -17952.5
$ws.Range("N65").Value = -38338.63499999999
$ws.Range("H86").Value = 4383.35
$ws.Range("I86").Value = 3680.2222
$ws.Range("K86").Value = 3680.2222
$ws.Range("M86").Value = -2557.2222
$ws.Range("H89").Value = 4383.35
$ws.Range("I89").Value = 3680.2222
$ws.Range("K89").Value = 18401.111
$ws.Range("M89").Value = -12785.111
$ws.Range("H99").Value = 909.17645
$ws.Range("I99").Value = 571.3333
$ws.Range("K99").Value = 1713.9999
$ws.Range("M99").Value = -215.9999
$ws.Range("H106").Value = 2013.5714
$ws.Range("I106").Value = 1848.1666
$ws.Range("K106").Value = 1848.1666
$ws.Range("M106").Value = -1217.1666
$ws.Range("H116").Value = 2675.1
$ws.Range("J116").Value = 3091.2
$ws.Range("L116").Value = 3091.2
$ws.Range("N116").Value = -9975.200000000001
$ws.Range("H131").Value = 4948.4287
$ws.Range("I131").Value = 2817.8
$ws.Range("K131").Value = 8453.400000000001
$ws.Range("M131").Value = -3413.400000000001
$ws.Range("H132").Value = 2402.926
$ws.Range("I132").Value = 1970
$ws.Range("K132").Value = 5910
$ws.Range("M132").Value = -3380
$ws.Range("H135").Value = 798.3333
$ws.Range("I135").Value = 499.8
$ws.Range("K135").Value = 4498.2
$ws.Range("M135").Value = -1963.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6672
$ws.Range("I32").Value = 4905.423
$ws.Range("K32").Value = 4905.423
$ws.Range("M32").Value = -4618.423
$ws.Range("H110").Value = 19121.65
$ws.Range("I110").Value = 21826.354
$ws.Range("K110").Value = 21826.354
$ws.Range("M110").Value = -19781.354
$ws.Range("H122").Value = 4104.7334
$ws.Range("I122").Value = 3238.6
$ws.Range("K122").Value = 9715.799999999999
$ws.Range("M122").Value = -7265.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 1140.4
$ws.Range("I7").Value = 351
$ws.Range("J7").Value = 1666.6666
$ws.Range("K7").Value = 351
$ws.Range("L7").Value = 1666.6666
$ws.Range("M7").Value = -238
$ws.Range("N7").Value = -1892.6666
$ws.Range("H11").Value = 723
$ws.Range("J11").Value = 723
$ws.Range("L11").Value = 723
$ws.Range("N11").Value = -1003
$ws.Range("H94").Value = 1258.5385
$ws.Range("I94").Value = 827
$ws.Range("J94").Value = 1949
$ws.Range("K94").Value = 827
$ws.Range("L94").Value = 1949
$ws.Range("M94").Value = -376
$ws.Range("N94").Value = -2851

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 162.81818
$ws.Range("I5").Value = 116
$ws.Range("J5").Value = 219
$ws.Range("K5").Value = 116
$ws.Range("L5").Value = 219
$ws.Range("M5").Value = -4
$ws.Range("N5").Value = -443
$ws.Range("H19").Value = 2244.6667
$ws.Range("I19").Value = 2806.3635
$ws.Range("K19").Value = 2806.3635
$ws.Range("M19").Value = -2636.3635
$ws.Range("H24").Value = 2244.6667
$ws.Range("I24").Value = 2806.3635
$ws.Range("K24").Value = 2806.3635
$ws.Range("M24").Value = -2636.3635
$ws.Range("H58").Value = 2180.0417
$ws.Range("I58").Value = 1426
$ws.Range("J58").Value = 3436.7778
$ws.Range("K58").Value = 1426
$ws.Range("L58").Value = 3436.7778
$ws.Range("M58").Value = -1223
$ws.Range("N58").Value = -3842.7778
$ws.Range("H105").Value = 17279.715
$ws.Range("I105").Value = 1999
$ws.Range("J105").Value = 28740.25
$ws.Range("K105").Value = 1999
$ws.Range("L105").Value = 28740.25
$ws.Range("M105").Value = -252
$ws.Range("N105").Value = -32234.25
$ws.Range("H132").Value = 5036.5386
$ws.Range("I132").Value = 4567.3887
$ws.Range("J132").Value = 10666.333
$ws.Range("K132").Value = 13702.1661
$ws.Range("L132").Value = 31998.999
$ws.Range("M132").Value = -11172.1661
$ws.Range("N132").Value = -37058.999
$ws.Range("H136").Value = 2180.0417
$ws.Range("I136").Value = 1426
$ws.Range("J136").Value = 3436.7778
$ws.Range("K136").Value = 4278
$ws.Range("L136").Value = 10310.3334
$ws.Range("M136").Value = -1728
$ws.Range("N136").Value = -15410.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 2374.6365
$ws.Range("I132").Value = 1677.8
$ws.Range("J132").Value = 2955.3333
$ws.Range("K132").Value = 15100.2
$ws.Range("L132").Value = 26597.9997
$ws.Range("M132").Value = -12570.2
$ws.Range("N132").Value = -31657.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 5248682.5
$ws.Range("I3").Value = 6560662
$ws.Range("J3").Value = 765
$ws.Range("K3").Value = 6560662
$ws.Range("L3").Value = 765
$ws.Range("M3").Value = -6560546
$ws.Range("N3").Value = -997
$ws.Range("H14").Value = 3992175
$ws.Range("I14").Value = 4213906.5
$ws.Range("J14").Value = 1005
$ws.Range("K14").Value = 4213906.5
$ws.Range("L14").Value = 1005
$ws.Range("M14").Value = -4213738.5
$ws.Range("N14").Value = -1341
$ws.Range("H97").Value = 3729.25
$ws.Range("I97").Value = 3616.0908
$ws.Range("J97").Value = 3978.2
$ws.Range("K97").Value = 3616.0908
$ws.Range("L97").Value = 3978.2
$ws.Range("M97").Value = -3120.0908
$ws.Range("N97").Value = -4970.2
$ws.Range("H102").Value = 2496.7693
$ws.Range("I102").Value = 1876.25
$ws.Range("J102").Value = 4565.1665
$ws.Range("K102").Value = 1876.25
$ws.Range("L102").Value = 4565.1665
$ws.Range("M102").Value = -254.25
$ws.Range("N102").Value = -7809.1665
$ws.Range("H113").Value = 3471.611
$ws.Range("I113").Value = 1967.3334
$ws.Range("J113").Value = 4223.75
$ws.Range("K113").Value = 1967.3334
$ws.Range("L113").Value = 4223.75
$ws.Range("M113").Value = 202.6666
$ws.Range("N113").Value = -8563.75
$ws.Range("H126").Value = 4635.222
$ws.Range("I126").Value = 4467.9165
$ws.Range("J126").Value = 4769.067
$ws.Range("K126").Value = 13403.7495
$ws.Range("L126").Value = 14307.201
$ws.Range("M126").Value = -10933.7495
$ws.Range("N126").Value = -19247.201
$ws.Range("H132").Value = 2667.05
$ws.Range("I132").Value = 2492.074
$ws.Range("J132").Value = 3030.4614
$ws.Range("K132").Value = 7476.222
$ws.Range("L132").Value = 9091.3842
$ws.Range("M132").Value = -4946.222
$ws.Range("N132").Value = -14151.3842

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3782.4666
$ws.Range("I7").Value = 3522.4546
$ws.Range("J7").Value = 4497.5
$ws.Range("K7").Value = 3522.4546
$ws.Range("L7").Value = 4497.5
$ws.Range("M7").Value = -3410.4546
$ws.Range("N7").Value = -4721.5
$ws.Range("H40").Value = 3575.9092
$ws.Range("I40").Value = 3715.611
$ws.Range("J40").Value = 2947.25
$ws.Range("K40").Value = 3715.611
$ws.Range("L40").Value = 2947.25
$ws.Range("M40").Value = -3579.611
$ws.Range("N40").Value = -3219.25
$ws.Range("H46").Value = 1172.5128
$ws.Range("J46").Value = 2891.2222
$ws.Range("L46").Value = 2891.2222
$ws.Range("N46").Value = -3267.2222
$ws.Range("H55").Value = 755.5294
$ws.Range("I55").Value = 404.2
$ws.Range("J55").Value = 1257.4286
$ws.Range("K55").Value = 404.2
$ws.Range("L55").Value = 1257.4286
$ws.Range("M55").Value = -231.2
$ws.Range("N55").Value = -1603.4286
$ws.Range("H100").Value = 3197.3333
$ws.Range("I100").Value = 2396.6667
$ws.Range("K100").Value = 2396.6667
$ws.Range("M100").Value = -1855.6667
$ws.Range("H122").Value = 4265.2
$ws.Range("I122").Value = 3353.3333
$ws.Range("K122").Value = 10059.9999
$ws.Range("M122").Value = -7609.999899999999
$ws.Range("H126").Value = 3782.4666
$ws.Range("I126").Value = 3522.4546
$ws.Range("J126").Value = 4497.5
$ws.Range("K126").Value = 10567.3638
$ws.Range("L126").Value = 13492.5
$ws.Range("M126").Value = -8097.363799999999
$ws.Range("N126").Value = -18432.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 50000
$ws.Range("J74").Value = 50000
$ws.Range("L74").Value = 50000
$ws.Range("N74").Value = -51872
$ws.Range("H77").Value = 50000
$ws.Range("J77").Value = 50000
$ws.Range("L77").Value = 150000
$ws.Range("N77").Value = -159360
$ws.Range("H95").Value = 30000
$ws.Range("J95").Value = 30000
$ws.Range("L95").Value = 30000
$ws.Range("N95").Value = -35492
$ws.Range("H96").Value = 5011.5
$ws.Range("I96").Value = 3824.5
$ws.Range("K96").Value = 3824.5
$ws.Range("M96").Value = -2451.5
$ws.Range("H126").Value = 7719.1
$ws.Range("I126").Value = 9199
$ws.Range("J126").Value = 1799.5
$ws.Range("K126").Value = 27597
$ws.Range("L126").Value = 5398.5
$ws.Range("M126").Value = -25127
$ws.Range("N126").Value = -10338.5
$ws.Range("H136").Value = 1938.4445
$ws.Range("I136").Value = 1910.3077
$ws.Range("K136").Value = 5730.9231
$ws.Range("M136").Value = -3180.9231
